# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (F) and "最低票价" (G) columns across all four
# sheets to match the newly scraped snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1394
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 13364
$ws.Range("G4").Value = 85
$ws.Range("F5").Value = 768
$ws.Range("F10").Value = 1914
$ws.Range("F11").Value = 49
$ws.Range("F13").Value = 20801
$ws.Range("G13").Value = 68
$ws.Range("F14").Value = 543
$ws.Range("F15").Value = 223
$ws.Range("F16").Value = 386
$ws.Range("F18").Value = 373
$ws.Range("F20").Value = 317
$ws.Range("F21").Value = 161
$ws.Range("F22").Value = 142
$ws.Range("F25").Value = 283
$ws.Range("F27").Value = 1359
$ws.Range("F28").Value = 61
$ws.Range("F29").Value = 378
$ws.Range("F30").Value = 78

# --- Sheet 2: 演出 -------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 300
$ws.Range("G3").Value = 199
$ws.Range("F4").Value = 4477
$ws.Range("F5").Value = 193
$ws.Range("F7").Value = 15
$ws.Range("F8").Value = 27
$ws.Range("F11").Value = 386
$ws.Range("F14").Value = 8
$ws.Range("F20").Value = 15

# --- Sheet 3: 本地生活 ---------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 906
$ws.Range("F3").Value = 4432
$ws.Range("F4").Value = 102

# --- Sheet 4: 全部类型 ---------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 906
$ws.Range("F5").Value = 1394
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 13364
$ws.Range("G6").Value = 85
$ws.Range("F7").Value = 300
$ws.Range("G7").Value = 199
$ws.Range("F8").Value = 768
$ws.Range("F9").Value = 4432
$ws.Range("F13").Value = 1914
$ws.Range("F14").Value = 49
$ws.Range("F16").Value = 102
$ws.Range("F17").Value = 20802
$ws.Range("G17").Value = 68
$ws.Range("F18").Value = 543
$ws.Range("F19").Value = 4477
$ws.Range("F20").Value = 223
$ws.Range("F21").Value = 193
$ws.Range("F22").Value = 193
$ws.Range("F23").Value = 386
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 27
$ws.Range("F29").Value = 386
$ws.Range("F30").Value = 373
$ws.Range("F32").Value = 317
$ws.Range("F33").Value = 161
$ws.Range("F34").Value = 142
$ws.Range("F38").Value = 8
$ws.Range("F40").Value = 283
$ws.Range("F42").Value = 1359
$ws.Range("F43").Value = 61
$ws.Range("F45").Value = 378
$ws.Range("F46").Value = 78
$ws.Range("F51").Value = 15
